# Updates cryptos.xlsx price (D) and volume (E) columns with fresh values.
# Values are plain text (matching the source sheet's inlineStr cells), so
# force text number-format before assignment to avoid Excel auto-coercing
# "307.65" / "2.68%" into numeric/percentage values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.65"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.68%"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.40%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.069"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.73%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08114"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.24%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.940"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.39%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.165"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.33%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.827"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.62%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9386"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.40%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1370"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-3.91%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1917"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.83%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09222"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.04%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03512"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.44%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09910"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.44%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001438"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2.91%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005867"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.63%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.626"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "3.03%"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.72%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3429"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.56%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1346"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "4.14%"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.97%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2533"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.21%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04403"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.33%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001236"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.65%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004765"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.24%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "5.53%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003131"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "4.35%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02023"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "6.65%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05106"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "8.49%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01126"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "16.08%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007668"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.21%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1378"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.64%"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.63%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01134"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "21.54%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006319"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.28%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.08%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.57"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-1.50%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001191"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-28.18%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.08%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.08%"

